# Updates Price (column D) and Volume(1h) (column E) cells to match the
# latest scraped crypto values. Price cells are written with a leading
# apostrophe (Excel's force-text marker) so number-looking strings such as
# "1.00" or "60.993.17" stay literal text instead of being parsed into
# numeric values; the Style reset afterwards keeps the cell on the default
# "Normal" style (no NumberFormat/quotePrefix left behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'" + '60.993.17'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.28%  '

$ws.Cells.Item(3, 4).Value = "'" + '3.381.90'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.10%  '

$ws.Cells.Item(4, 4).Value = "'" + '1.00'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$ws.Cells.Item(5, 4).Value = "'" + '571.20'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.61%  '

$ws.Cells.Item(6, 4).Value = "'" + '140.81'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.65%  '

$ws.Cells.Item(8, 5).Value = '  +0.27%  '

$ws.Cells.Item(9, 4).Value = "'" + '7.67'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.75%  '

$ws.Cells.Item(10, 5).Value = '  -0.58%  '

$ws.Cells.Item(11, 4).Value = "'" + '0.388'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.17%  '

$ws.Cells.Item(12, 4).Value = "'" + '3.960.54'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.13%  '

$ws.Cells.Item(13, 5).Value = '  +2.08%  '

$ws.Cells.Item(14, 4).Value = "'" + '27.86'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.50%  '

$ws.Cells.Item(15, 4).Value = "'" + '3.377.95'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.23%  '

$ws.Cells.Item(16, 5).Value = '  +0.37%  '

$ws.Cells.Item(17, 4).Value = "'" + '61.099.54'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.19%  '

$ws.Cells.Item(18, 4).Value = "'" + '6.11'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.87%  '

$ws.Cells.Item(19, 4).Value = "'" + '13.60'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.67%  '

$ws.Cells.Item(20, 4).Value = "'" + '8.89'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -1.75%  '

$ws.Cells.Item(21, 4).Value = "'" + '383.29'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.83%  '

$ws.Cells.Item(22, 4).Value = "'" + '75.85'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +3.47%  '

$ws.Cells.Item(23, 5).Value = '  -0.96%  '

$ws.Cells.Item(24, 5).Value = '  +0.00%  '

$ws.Cells.Item(25, 5).Value = '  -0.87%  '

$ws.Cells.Item(26, 4).Value = "'" + '0.186'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +4.07%  '

$ws.Cells.Item(27, 5).Value = '  +0.01%  '

$ws.Cells.Item(28, 5).Value = '  -1.87%  '

$ws.Cells.Item(29, 4).Value = "'" + '7.95'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.36%  '

$ws.Cells.Item(30, 4).Value = "'" + '2.13'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.00%  '

$ws.Cells.Item(31, 5).Value = '  -0.04%  '

$ws.Cells.Item(32, 4).Value = "'" + '1.36'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.61%  '

$ws.Cells.Item(33, 4).Value = "'" + '23.34'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.29%  '

$ws.Cells.Item(34, 5).Value = '  +0.59%  '

$ws.Cells.Item(35, 4).Value = "'" + '165.64'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.43%  '

$ws.Cells.Item(36, 4).Value = "'" + '3.417.54'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

$ws.Cells.Item(37, 4).Value = "'" + '4.98'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.54%  '

$ws.Cells.Item(38, 5).Value = '  -2.26%  '

$ws.Cells.Item(39, 4).Value = "'" + '0.0766'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.68%  '

$ws.Cells.Item(40, 4).Value = "'" + '26.42'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.15%  '

$ws.Cells.Item(41, 5).Value = '  +0.00%  '

$ws.Cells.Item(42, 4).Value = "'" + '0.778'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.66%  '

$ws.Cells.Item(43, 5).Value = '  -1.70%  '

$ws.Cells.Item(44, 5).Value = '  -2.34%  '

$ws.Cells.Item(45, 5).Value = '  +0.58%  '

$ws.Cells.Item(46, 4).Value = "'" + '2.458.37'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.48%  '

$ws.Cells.Item(47, 4).Value = "'" + '22.84'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.45%  '

$ws.Cells.Item(48, 4).Value = "'" + '6.63'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -2.29%  '

$ws.Cells.Item(49, 5).Value = '  +10.43%  '

$ws.Cells.Item(50, 4).Value = "'" + '0.0261'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.95%  '

$ws.Cells.Item(51, 5).Value = '  -1.47%  '
